{"js": "// The document's Heading1/2/3 paragraphs each carry a pair of Word TOC\n// anchor bookmarks: an older \"_Toc24446xxx\" one (left untouched) and a\n// newer \"_Toc24731xxx\" one. This edit corresponds to Word regenerating the\n// table of contents (\"Update entire table\" / re-save), which reassigns a\n// fresh set of \"_Toc\" bookmark names (here \"_Toc25148430\"..\"_Toc25148443\")\n// to the same headings, in the same order, while keeping their bookmark\n// ids stable.\n//\n// Office.js has no \"rename bookmark\" primitive, so each rename is done by\n// capturing the existing bookmark's range, inserting a new bookmark with\n// the new name over that same range, and then deleting the old bookmark.\n\nconst renames = [\n  [\"_Toc24731331\", \"_Toc25148430\"],\n  [\"_Toc24731332\", \"_Toc25148431\"],\n  [\"_Toc24731333\", \"_Toc25148432\"],\n  [\"_Toc24731334\", \"_Toc25148433\"],\n  [\"_Toc24731335\", \"_Toc25148434\"],\n  [\"_Toc24731336\", \"_Toc25148435\"],\n  [\"_Toc24731337\", \"_Toc25148436\"],\n  [\"_Toc24731338\", \"_Toc25148437\"],\n  [\"_Toc24731339\", \"_Toc25148438\"],\n  [\"_Toc24731340\", \"_Toc25148439\"],\n  [\"_Toc24731341\", \"_Toc25148440\"],\n  [\"_Toc24731342\", \"_Toc25148441\"],\n  [\"_Toc24731343\", \"_Toc25148442\"],\n  [\"_Toc24731344\", \"_Toc25148443\"],\n];\n\nfor (const [oldName, newName] of renames) {\n  const exists = context.document.bookmarks.exists(oldName);\n  await context.sync();\n  if (!exists.value) {\n    continue;\n  }\n\n  const bookmark = context.document.bookmarks.getByName(oldName);\n  const range = bookmark.getRange();\n  // Make sure the range is resolved before we mutate bookmarks on it.\n  await context.sync();\n\n  range.insertBookmark(newName);\n  context.document.deleteBookmark(oldName);\n  await context.sync();\n}\n", "ps1": "# The document's Heading1/2/3 paragraphs each carry a pair of Word TOC\n# anchor bookmarks: an older \"_Toc24446xxx\" one (left untouched) and a\n# newer \"_Toc24731xxx\" one. This edit corresponds to Word regenerating the\n# table of contents (\"Update entire table\" / re-save), which reassigns a\n# fresh set of \"_Toc\" bookmark names (here \"_Toc25148430\"..\"_Toc25148443\")\n# to the same headings, in the same order, while keeping the bookmark ids\n# stable.\n#\n# The Word object model has no \"rename bookmark\" call, so each rename is\n# done by adding a new bookmark with the new name over the old bookmark's\n# Range, then deleting the old bookmark.\n\n$d = $word.ActiveDocument\n\n$renames = @(\n    @('_Toc24731331', '_Toc25148430'),\n    @('_Toc24731332', '_Toc25148431'),\n    @('_Toc24731333', '_Toc25148432'),\n    @('_Toc24731334', '_Toc25148433'),\n    @('_Toc24731335', '_Toc25148434'),\n    @('_Toc24731336', '_Toc25148435'),\n    @('_Toc24731337', '_Toc25148436'),\n    @('_Toc24731338', '_Toc25148437'),\n    @('_Toc24731339', '_Toc25148438'),\n    @('_Toc24731340', '_Toc25148439'),\n    @('_Toc24731341', '_Toc25148440'),\n    @('_Toc24731342', '_Toc25148441'),\n    @('_Toc24731343', '_Toc25148442'),\n    @('_Toc24731344', '_Toc25148443')\n)\n\nforeach ($pair in $renames) {\n    $oldName = $pair[0]\n    $newName = $pair[1]\n\n    if ($d.Bookmarks.Exists($oldName)) {\n        $rng = $d.Bookmarks.Item($oldName).Range\n        $d.Bookmarks.Add($newName, $rng)\n        $d.Bookmarks.Item($oldName).Delete()\n    }\n}\n"}
